$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.010.48'
$ws.Range('E2').Value = '  -2.05%  '
$ws.Range('D3').Value = '3.418.11'
$ws.Range('E3').Value = '  -1.52%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '578.15'
$ws.Range('D6').Value = '153.10'
$ws.Range('E6').Value = '  +3.68%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +1.31%  '
$ws.Range('D9').Value = '8.03'
$ws.Range('E9').Value = '  +3.72%  '
$ws.Range('E10').Value = '  -0.69%  '
$ws.Range('D11').Value = '0.417'
$ws.Range('E11').Value = '  +3.13%  '
$ws.Range('D12').Value = '4.004.21'
$ws.Range('E12').Value = '  -1.53%  '
$ws.Range('E13').Value = '  +0.59%  '
$ws.Range('D14').Value = '28.70'
$ws.Range('E14').Value = '  -2.81%  '
$ws.Range('D15').Value = '3.421.90'
$ws.Range('E15').Value = '  -1.64%  '
$ws.Range('E16').Value = '  -0.40%  '
$ws.Range('D17').Value = '62.045.96'
$ws.Range('E17').Value = '  -2.06%  '
$ws.Range('D18').Value = '6.50'
$ws.Range('E18').Value = '  +1.76%  '
$ws.Range('D19').Value = '14.50'
$ws.Range('E19').Value = '  +0.07%  '
$ws.Range('E20').Value = '  -4.21%  '
$ws.Range('D21').Value = '382.64'
$ws.Range('E21').Value = '  -1.67%  '
$ws.Range('E22').Value = '  +1.00%  '
$ws.Range('D23').Value = '75.26'
$ws.Range('E23').Value = '  +1.00%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').Value = '3.560.26'
$ws.Range('E25').Value = '  -1.51%  '
$ws.Range('E26').Value = '  -3.56%  '
$ws.Range('E27').Value = '  -1.51%  '
$ws.Range('D28').Value = '7.67'
$ws.Range('E28').Value = '  +0.31%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('D30').Value = '7.92'
$ws.Range('E30').Value = '  -3.63%  '
$ws.Range('E31').Value = '  -1.19%  '
$ws.Range('E32').Value = '  -0.07%  '
$ws.Range('E33').Value = '  -0.70%  '
$ws.Range('D34').Value = '23.25'
$ws.Range('E34').Value = '  -1.09%  '
$ws.Range('D35').Value = '5.52'
$ws.Range('E35').Value = '  +3.49%  '
$ws.Range('E36').Value = '  -0.55%  '
$ws.Range('E37').Value = '  -2.64%  '
$ws.Range('D38').Value = '168.23'
$ws.Range('E38').Value = '  -0.29%  '
$ws.Range('D39').Value = '30.89'
$ws.Range('E39').Value = '  -3.70%  '
$ws.Range('D40').Value = '3.453.70'
$ws.Range('E40').Value = '  -1.64%  '
$ws.Range('E41').Value = '  +2.52%  '
$ws.Range('D42').Value = '42.67'
$ws.Range('E42').Value = '  +0.60%  '
$ws.Range('D43').Value = '0.780'
$ws.Range('E43').Value = '  -1.98%  '
$ws.Range('D44').Value = '4.42'
$ws.Range('E44').Value = '  +0.61%  '
$ws.Range('E45').Value = '  -3.85%  '
$ws.Range('E46').Value = '  -3.41%  '
$ws.Range('D47').Value = '2.551.16'
$ws.Range('E47').Value = '  -1.75%  '
$ws.Range('E48').Value = '  +0.52%  '
$ws.Range('D49').Value = '22.63'
$ws.Range('E49').Value = '  -1.80%  '
$ws.Range('D50').Value = '2.19'
$ws.Range('E50').Value = '  -5.40%  '
$ws.Range('E51').Value = '  -0.07%  '
